# 20150113 +++++++ cs-厂商 end
# Update the "搜索列表" (search list) rows in the cs-厂商 comparison table:
#   - Row 10, column A gains the extra "搜索列表 - 排序.psd" psd reference
#   - A brand new row 11 is appended for "搜索列表.psd," -> "search-list.html"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new shared string for row 11 col A first, then update row 10 col A,
# so the shared-string table ends up in the same order as the saved workbook.
$ws.Range("A11").Value = "搜索列表.psd,"
$ws.Range("A10").Value = "搜索列表 - 大图版面.psd,搜索列表 - 排序.psd"
$ws.Range("B10").Value = "search-list-b.html"
$ws.Range("B11").Value = "search-list.html"

# Match the workbook's saved selection state
$ws.Cells.Item(10, 3).Select()
